# Regenerate the "K" column (column G) values in the save_data sheet.
# The workbook records per-game pitching stats; column G (header "K") is
# recomputed from the underlying strikeout data and rewritten here to the
# freshly-derived counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeout) values, keyed by data row (row 2 = first data record).
$kValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 0
    6  = 1
    7  = 1
    8  = 0
    9  = 1
    10 = 2
    11 = 1
    12 = 1
    13 = 0
    14 = 0
    15 = 1
    16 = 2
    17 = 1
    18 = 2
    19 = 1
    20 = 1
    22 = 2
    23 = 0
    24 = 0
    25 = 2
    26 = 2
    27 = 2
    28 = 1
    29 = 0
    30 = 1
    31 = 2
    32 = 1
    33 = 1
    34 = 1
    35 = 4
    36 = 1
    37 = 1
    38 = 2
    40 = 1
    41 = 1
    42 = 2
    43 = 2
    44 = 3
    45 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
